$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(2)

# tx9 (GroupItems index 7, id=10)
$gi = $grp.GroupItems.Item(7)
$gi.Left = 441.0425266850394
$gi.Top = 375.8163919527559

# tx10 (GroupItems index 8, id=11)
$gi = $grp.GroupItems.Item(8)
$gi.Left = 464.2785189370079
$gi.Top = 399.94111636220475

# tx11 (GroupItems index 9, id=12)
$gi = $grp.GroupItems.Item(9)
$gi.Left = 246.84865641732284
$gi.Top = 284.29426596850396

# tx12 (GroupItems index 10, id=13)
$gi = $grp.GroupItems.Item(10)
$gi.Left = 293.33308486614175
$gi.Top = 308.41899137795275

# tx13 (GroupItems index 11, id=14)
$gi = $grp.GroupItems.Item(11)
$gi.Left = 338.15739457480316
$gi.Top = 219.82062592125985

# tx14 (GroupItems index 12, id=15)
$gi = $grp.GroupItems.Item(12)
$gi.Left = 351.91142332283465
$gi.Top = 247.13787101574803

# tx15 (GroupItems index 13, id=16)
$gi = $grp.GroupItems.Item(13)
$gi.Left = 374.10072366141736
$gi.Top = 175.63079140157478

# tx16 (GroupItems index 14, id=17)
$gi = $grp.GroupItems.Item(14)
$gi.Left = 398.2869421338583
$gi.Top = 203.14811023622048
